# Auto-generated: updates Asura_Profits-style numeric cells (H:N) across all leve sheets.
# Source: scheduled runner recomputing currentAveragePrice* / LevePrice* / LeveProfit* columns.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1990
$ws.Range("I40").Value = 1480
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 1480
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -1305
$ws.Range("N40").Value = -2850
$ws.Range("H64").Value = 4203.8535
$ws.Range("I64").Value = 3427.2856
$ws.Range("J64").Value = 4363.7354
$ws.Range("K64").Value = 3427.2856
$ws.Range("L64").Value = 4363.7354
$ws.Range("M64").Value = -3179.2856
$ws.Range("N64").Value = -4859.7354
$ws.Range("H67").Value = 4203.8535
$ws.Range("I67").Value = 3427.2856
$ws.Range("J67").Value = 4363.7354
$ws.Range("K67").Value = 3427.2856
$ws.Range("L67").Value = 4363.7354
$ws.Range("M67").Value = -2569.2856
$ws.Range("N67").Value = -6079.7354
$ws.Range("H76").Value = 4387.5
$ws.Range("I76").Value = 4800
$ws.Range("J76").Value = 4328.5713
$ws.Range("K76").Value = 4800
$ws.Range("L76").Value = 4328.5713
$ws.Range("M76").Value = -4485
$ws.Range("N76").Value = -4958.5713
$ws.Range("H79").Value = 4387.5
$ws.Range("I79").Value = 4800
$ws.Range("J79").Value = 4328.5713
$ws.Range("K79").Value = 4800
$ws.Range("L79").Value = 4328.5713
$ws.Range("M79").Value = -3708
$ws.Range("N79").Value = -6512.5713
$ws.Range("H98").Value = 6828.68
$ws.Range("I98").Value = 6125.625
$ws.Range("J98").Value = 8078.5557
$ws.Range("K98").Value = 6125.625
$ws.Range("L98").Value = 8078.5557
$ws.Range("M98").Value = -4627.625
$ws.Range("N98").Value = -11074.5557
$ws.Range("H122").Value = 6828.68
$ws.Range("I122").Value = 6125.625
$ws.Range("J122").Value = 8078.5557
$ws.Range("K122").Value = 18376.875
$ws.Range("L122").Value = 24235.6671
$ws.Range("M122").Value = -15926.875
$ws.Range("N122").Value = -29135.6671
$ws.Range("H129").Value = 947.30646
$ws.Range("J129").Value = 1105.9375
$ws.Range("L129").Value = 3317.8125
$ws.Range("N129").Value = -13317.8125
$ws.Range("H137").Value = 1454.8431
$ws.Range("I137").Value = 1241.0256
$ws.Range("K137").Value = 3723.0768
$ws.Range("M137").Value = -1173.0768

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 704.8
$ws.Range("I45").Value = 606
$ws.Range("K45").Value = 606
$ws.Range("M45").Value = -229
$ws.Range("H61").Value = 2368.52
$ws.Range("I61").Value = 2246.1428
$ws.Range("J61").Value = 3011
$ws.Range("K61").Value = 2246.1428
$ws.Range("L61").Value = 3011
$ws.Range("M61").Value = -2034.1428
$ws.Range("N61").Value = -3435
$ws.Range("H74").Value = 1642.9032
$ws.Range("I74").Value = 1496.3334
$ws.Range("K74").Value = 1496.3334
$ws.Range("M74").Value = -622.3334
$ws.Range("H77").Value = 1642.9032
$ws.Range("I77").Value = 1496.3334
$ws.Range("K77").Value = 7481.666999999999
$ws.Range("M77").Value = -3113.666999999999
$ws.Range("H88").Value = 2850
$ws.Range("I88").Value = 2299
$ws.Range("J88").Value = 3401
$ws.Range("K88").Value = 2299
$ws.Range("L88").Value = 3401
$ws.Range("M88").Value = -1893
$ws.Range("N88").Value = -4213
$ws.Range("H91").Value = 2850
$ws.Range("I91").Value = 2299
$ws.Range("J91").Value = 3401
$ws.Range("K91").Value = 2299
$ws.Range("L91").Value = 3401
$ws.Range("M91").Value = -895
$ws.Range("N91").Value = -6209
$ws.Range("H122").Value = 2621.48
$ws.Range("I122").Value = 2641.9546
$ws.Range("J122").Value = 2471.3333
$ws.Range("K122").Value = 7925.8638
$ws.Range("L122").Value = 7413.999899999999
$ws.Range("M122").Value = -5475.8638
$ws.Range("N122").Value = -12313.9999
$ws.Range("H132").Value = 393560.1
$ws.Range("I132").Value = 466072.66
$ws.Range("K132").Value = 1398217.98
$ws.Range("M132").Value = -1395687.98
$ws.Range("H136").Value = 2368.52
$ws.Range("I136").Value = 2246.1428
$ws.Range("J136").Value = 3011
$ws.Range("K136").Value = 6738.428400000001
$ws.Range("L136").Value = 9033
$ws.Range("M136").Value = -4188.428400000001
$ws.Range("N136").Value = -14133

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 38000
$ws.Range("J88").Value = 38000
$ws.Range("L88").Value = 38000
$ws.Range("N88").Value = -38812
$ws.Range("H91").Value = 38000
$ws.Range("J91").Value = 38000
$ws.Range("L91").Value = 38000
$ws.Range("N91").Value = -40808
$ws.Range("H105").Value = 3638.7144
$ws.Range("I105").Value = 3161.8333
$ws.Range("J105").Value = 6500
$ws.Range("K105").Value = 3161.8333
$ws.Range("L105").Value = 6500
$ws.Range("M105").Value = -1414.8333
$ws.Range("N105").Value = -9994
$ws.Range("H134").Value = 296184.44
$ws.Range("I134").Value = 372008.4
$ws.Range("J134").Value = 3720.5715
$ws.Range("K134").Value = 1116025.2
$ws.Range("L134").Value = 11161.7145
$ws.Range("M134").Value = -1113490.2
$ws.Range("N134").Value = -16231.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2791.36
$ws.Range("I31").Value = 1690.7858
$ws.Range("J31").Value = 4192.091
$ws.Range("K31").Value = 1690.7858
$ws.Range("L31").Value = 4192.091
$ws.Range("M31").Value = -1395.7858
$ws.Range("N31").Value = -4782.091
$ws.Range("H34").Value = 2791.36
$ws.Range("I34").Value = 1690.7858
$ws.Range("J34").Value = 4192.091
$ws.Range("K34").Value = 1690.7858
$ws.Range("L34").Value = 4192.091
$ws.Range("M34").Value = -1488.7858
$ws.Range("N34").Value = -4596.091
$ws.Range("H58").Value = 1324062.8
$ws.Range("I58").Value = 1765011.2
$ws.Range("K58").Value = 1765011.2
$ws.Range("M58").Value = -1764808.2
$ws.Range("H122").Value = 1337.3334
$ws.Range("I122").Value = 1046
$ws.Range("J122").Value = 1920
$ws.Range("K122").Value = 3138
$ws.Range("L122").Value = 5760
$ws.Range("M122").Value = -688
$ws.Range("N122").Value = -10660
$ws.Range("H132").Value = 339707.06
$ws.Range("I132").Value = 484240.75
$ws.Range("K132").Value = 1452722.25
$ws.Range("M132").Value = -1450192.25
$ws.Range("H134").Value = 1280.2424
$ws.Range("I134").Value = 1141.5555
$ws.Range("K134").Value = 3424.6665
$ws.Range("M134").Value = -889.6664999999998
$ws.Range("H136").Value = 1324062.8
$ws.Range("I136").Value = 1765011.2
$ws.Range("K136").Value = 5295033.6
$ws.Range("M136").Value = -5292483.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4746.625
$ws.Range("J39").Value = 4746.625
$ws.Range("L39").Value = 14239.875
$ws.Range("N39").Value = -14827.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6272.2383
$ws.Range("I70").Value = 5251
$ws.Range("J70").Value = 6900.6924
$ws.Range("K70").Value = 5251
$ws.Range("L70").Value = 6900.6924
$ws.Range("M70").Value = -4981
$ws.Range("N70").Value = -7440.6924
$ws.Range("H73").Value = 6272.2383
$ws.Range("I73").Value = 5251
$ws.Range("J73").Value = 6900.6924
$ws.Range("K73").Value = 5251
$ws.Range("L73").Value = 6900.6924
$ws.Range("M73").Value = -4315
$ws.Range("N73").Value = -8772.6924
$ws.Range("H122").Value = 3926.5715
$ws.Range("I122").Value = 3591.1562
$ws.Range("J122").Value = 4999.9
$ws.Range("K122").Value = 10773.4686
$ws.Range("L122").Value = 14999.7
$ws.Range("M122").Value = -8323.4686
$ws.Range("N122").Value = -19899.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9732.105
$ws.Range("I2").Value = 5000
$ws.Range("K2").Value = 5000
$ws.Range("M2").Value = -4888
$ws.Range("H7").Value = 3189.8
$ws.Range("I7").Value = 3314.8333
$ws.Range("J7").Value = 3002.25
$ws.Range("K7").Value = 3314.8333
$ws.Range("L7").Value = 3002.25
$ws.Range("M7").Value = -3202.8333
$ws.Range("N7").Value = -3226.25
$ws.Range("H61").Value = 11540.429
$ws.Range("I61").Value = 15962.429
$ws.Range("J61").Value = 2696.4285
$ws.Range("K61").Value = 15962.429
$ws.Range("L61").Value = 2696.4285
$ws.Range("M61").Value = -15760.429
$ws.Range("N61").Value = -3100.4285
$ws.Range("H113").Value = 11540.429
$ws.Range("I113").Value = 15962.429
$ws.Range("J113").Value = 2696.4285
$ws.Range("K113").Value = 15962.429
$ws.Range("L113").Value = 2696.4285
$ws.Range("M113").Value = -13792.429
$ws.Range("N113").Value = -7036.4285
$ws.Range("H126").Value = 3189.8
$ws.Range("I126").Value = 3314.8333
$ws.Range("J126").Value = 3002.25
$ws.Range("K126").Value = 9944.499899999999
$ws.Range("L126").Value = 9006.75
$ws.Range("M126").Value = -7474.499899999999
$ws.Range("N126").Value = -13946.75
$ws.Range("H136").Value = 3576.6785
$ws.Range("I136").Value = 3699.8635
$ws.Range("J136").Value = 3125
$ws.Range("K136").Value = 11099.5905
$ws.Range("L136").Value = 9375
$ws.Range("M136").Value = -8549.5905
$ws.Range("N136").Value = -14475

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 16600.334
$ws.Range("I21").Value = 50000
$ws.Range("J21").Value = 9920.4
$ws.Range("K21").Value = 50000
$ws.Range("L21").Value = 9920.4
$ws.Range("M21").Value = -49765
$ws.Range("N21").Value = -10390.4
$ws.Range("H35").Value = 16600.334
$ws.Range("I35").Value = 50000
$ws.Range("J35").Value = 9920.4
$ws.Range("K35").Value = 50000
$ws.Range("L35").Value = 9920.4
$ws.Range("M35").Value = -49710
$ws.Range("N35").Value = -10500.4
$ws.Range("H113").Value = 1111.8
$ws.Range("I113").Value = 1111.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3335.4
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1165.4
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 19231970
$ws.Range("I122").Value = 25000880
$ws.Range("J122").Value = 2266.6667
$ws.Range("K122").Value = 75002640
$ws.Range("L122").Value = 6800.000100000001
$ws.Range("M122").Value = -75000190
$ws.Range("N122").Value = -11700.0001
$ws.Range("H126").Value = 2756
$ws.Range("I126").Value = 3012.6316
$ws.Range("J126").Value = 1780.8
$ws.Range("K126").Value = 9037.8948
$ws.Range("L126").Value = 5342.4
$ws.Range("M126").Value = -6567.8948
$ws.Range("N126").Value = -10282.4
$ws.Range("H132").Value = 2334.1155
$ws.Range("I132").Value = 2114.5
$ws.Range("J132").Value = 2590.3333
$ws.Range("K132").Value = 6343.5
$ws.Range("L132").Value = 7770.999899999999
$ws.Range("M132").Value = -3813.5
$ws.Range("N132").Value = -12830.9999
$ws.Range("H136").Value = 2681.2144
$ws.Range("I136").Value = 2926.0557
$ws.Range("J136").Value = 2240.5
$ws.Range("K136").Value = 8778.167099999999
$ws.Range("L136").Value = 6721.5
$ws.Range("M136").Value = -6228.167099999999
$ws.Range("N136").Value = -11821.5
